# Fix the property_category column ("building sheet" cells used to hold
# the wrong value "land" copy/pasted from the 土地 sheet; same mistake on
# the 汽車 (car) sheet). Correct each to match the sheet it lives on.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet — column I is property_category; rows 2-5 currently
# say "land", should say "building".
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"
$wsBuilding.Range("I4").Value = "building"
$wsBuilding.Range("I5").Value = "building"

# 汽車 (car) sheet — column H is property_category; row 2 currently says
# "land", should say "car".
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
